# QCA-AID-Codebook.xlsx CONFIG sheet update
# - lowercases the config keys
# - updates several config values
# - replaces the old ATTRIBUTE_LABELS[attribut3] / EXPORT_ANNOTATED_PDFS /
#   PDF_ANNOTATION_FUZZY_THRESHOLD / CODER_SETTINGS block with the new
#   coder_settings / manual_coding_enabled / export_annotated_pdfs /
#   pdf_annotation_fuzzy_threshold block
# - removes the trailing VALIDATION block (rows 23-38)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CONFIG")

# Remove the old VALIDATION rows (23-38) first so the dimension shrinks to A1:E22
$ws.Range("A23:E38").EntireRow.Delete()

# Row 2: MODEL_PROVIDER -> model_provider / OpenAI -> openai
$ws.Range("A2").Value = "model_provider"
$ws.Range("D2").Value = "openai"

# Row 3: MODEL_NAME -> model_name / gpt-4o-mini -> gpt-5-nano
$ws.Range("A3").Value = "model_name"
$ws.Range("D3").Value = "gpt-5-nano"

# Row 4: DATA_DIR -> data_dir
$ws.Range("A4").Value = "data_dir"

# Row 5: OUTPUT_DIR -> output_dir
$ws.Range("A5").Value = "output_dir"

# Row 6: CHUNK_SIZE -> chunk_size / 1000 -> 1200
$ws.Range("A6").Value = "chunk_size"
$ws.Range("D6").Value = 1200

# Row 7: CHUNK_OVERLAP -> chunk_overlap / 40 -> 50
$ws.Range("A7").Value = "chunk_overlap"
$ws.Range("D7").Value = 50

# Row 8: BATCH_SIZE -> batch_size / 5 -> 8
$ws.Range("A8").Value = "batch_size"
$ws.Range("D8").Value = 8

# Row 9: CODE_WITH_CONTEXT -> code_with_context
$ws.Range("A9").Value = "code_with_context"

# Row 10: MULTIPLE_CODINGS -> multiple_codings
$ws.Range("A10").Value = "multiple_codings"

# Row 11: MULTIPLE_CODING_THRESHOLD -> multiple_coding_threshold
$ws.Range("A11").Value = "multiple_coding_threshold"

# Row 12: ANALYSIS_MODE -> analysis_mode
$ws.Range("A12").Value = "analysis_mode"

# Row 13: REVIEW_MODE -> review_mode / consensus -> auto
$ws.Range("A13").Value = "review_mode"
$ws.Range("D13").Value = "auto"

# Row 14: ATTRIBUTE_LABELS -> attribute_labels / Hochschulprofil -> Fall
$ws.Range("A14").Value = "attribute_labels"
$ws.Range("D14").Value = "Fall"

# Row 15: ATTRIBUTE_LABELS -> attribute_labels / Akteur -> Typ
$ws.Range("A15").Value = "attribute_labels"
$ws.Range("D15").Value = "Typ"

# Row 16: was ATTRIBUTE_LABELS/attribut3 -> now coder_settings/[0]/temperature/0.3
$ws.Range("A16").Value = "coder_settings"
$ws.Range("B16").Value = "[0]"
$ws.Range("C16").Value = "temperature"
$ws.Range("D16").Value = 0.3

# Row 17: was EXPORT_ANNOTATED_PDFS -> now coder_settings/[0]/coder_id/auto_1
$ws.Range("A17").Value = "coder_settings"
$ws.Range("B17").Value = "[0]"
$ws.Range("C17").Value = "coder_id"
$ws.Range("D17").Value = "auto_1"

# Row 18: was PDF_ANNOTATION_FUZZY_THRESHOLD -> now coder_settings/[1]/temperature/0.5
$ws.Range("A18").Value = "coder_settings"
$ws.Range("B18").Value = "[1]"
$ws.Range("C18").Value = "temperature"
$ws.Range("D18").Value = 0.5

# Row 19: was CODER_SETTINGS/[0]/temperature/0.3 -> now coder_settings/[1]/coder_id/auto_2
$ws.Range("A19").Value = "coder_settings"
$ws.Range("B19").Value = "[1]"
$ws.Range("C19").Value = "coder_id"
$ws.Range("D19").Value = "auto_2"

# Row 20: was CODER_SETTINGS/[0]/coder_id/auto_1 -> now manual_coding_enabled / 0
$ws.Range("A20").Value = "manual_coding_enabled"
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = ""
$ws.Range("D20").Value = 0

# Row 21: was CODER_SETTINGS/[1]/temperature/0.5 -> now export_annotated_pdfs / 1
$ws.Range("A21").Value = "export_annotated_pdfs"
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = ""
$ws.Range("D21").Value = 1

# Row 22: was CODER_SETTINGS/[1]/coder_id/auto_2 -> now pdf_annotation_fuzzy_threshold / 0.85
$ws.Range("A22").Value = "pdf_annotation_fuzzy_threshold"
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = ""
$ws.Range("D22").Value = 0.85

Write-Host "CONFIG sheet updated. UsedRange:" $ws.UsedRange.Address()
